$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.872.37'
$ws.Range("E2").Value = '  -0.40%  '

$ws.Range("D3").Value = '2.360.92'
$ws.Range("E3").Value = '  -3.80%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.58'
$ws.Range("E5").Value = '  -1.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.60'
$ws.Range("E6").Value = '  -5.68%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.524'
$ws.Range("E8").Value = '  -10.36%  '

$ws.Range("D9").Value = '2.359.58'
$ws.Range("E9").Value = '  -3.83%  '

$ws.Range("E10").Value = '  -1.38%  '

$ws.Range("E11").Value = '  +0.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.23'
$ws.Range("E12").Value = '  -3.71%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.341'
$ws.Range("E13").Value = '  -2.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.79'
$ws.Range("E14").Value = '  -4.46%  '

$ws.Range("D15").Value = '2.790.31'
$ws.Range("E15").Value = '  -3.72%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '60.640.00'
$ws.Range("E16").Value = '  -0.68%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000161'
$ws.Range("E17").Value = '  -2.64%  '

$ws.Range("D18").Value = '2.362.71'
$ws.Range("E18").Value = '  -3.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.60'
$ws.Range("E19").Value = '  -4.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '315.19'
$ws.Range("E20").Value = '  -0.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.05'
$ws.Range("E21").Value = '  -2.30%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.58'
$ws.Range("E22").Value = '  -6.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.90'
$ws.Range("E24").Value = '  +2.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.25'
$ws.Range("E25").Value = '  -0.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.45'
$ws.Range("E26").Value = '  +12.03%  '

$ws.Range("E27").Value = '  -0.12%  '

$ws.Range("D28").Value = '2.483.06'
$ws.Range("E28").Value = '  -3.68%  '

$ws.Range("D29").Value = '0.0₃0898'
$ws.Range("E29").Value = '  -6.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.95'
$ws.Range("E30").Value = '  -2.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '503.20'
$ws.Range("E31").Value = '  -8.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.38'
$ws.Range("E32").Value = '  -4.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.145'
$ws.Range("E33").Value = '  -1.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.78'
$ws.Range("E34").Value = '  -5.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  -2.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.55'
$ws.Range("E37").Value = '  -4.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.56'
$ws.Range("E38").Value = '  +0.97%  '

$ws.Range("E39").Value = '  -1.19%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.25'
$ws.Range("E40").Value = '  -9.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.77'
$ws.Range("E41").Value = '  +0.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.34'
$ws.Range("E42").Value = '  -1.82%  '

$ws.Range("E43").Value = '  -0.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.19'
$ws.Range("E44").Value = '  -0.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.11'
$ws.Range("E45").Value = '  -7.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '138.42'
$ws.Range("E46").Value = '  -5.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.51'
$ws.Range("E47").Value = '  -1.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0511'
$ws.Range("E48").Value = '  -4.29%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.44'
$ws.Range("E49").Value = '  -8.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.570'
$ws.Range("E50").Value = '  -2.79%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0896'
$ws.Range("E51").Value = '  -3.57%  '
